# Update UnitMass ("C" column) example data values for the loading tables
# on the active worksheet (LoadingPC5.xlsx), as part of refreshing the
# example report output with data from the catalog system.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "+ loading" table (rows 2-21)
$ws.Range("C2").Value = 23
$ws.Range("C3").Value = 43
$ws.Range("C4").Value = 41
$ws.Range("C5").Value = 71
$ws.Range("C6").Value = 55
$ws.Range("C7").Value = 125
$ws.Range("C8").Value = 83
$ws.Range("C9").Value = 46
$ws.Range("C10").Value = 29
$ws.Range("C11").Value = 24
$ws.Range("C12").Value = 69
$ws.Range("C13").Value = 102
$ws.Range("C14").Value = 53
$ws.Range("C15").Value = 25
$ws.Range("C16").Value = 7
$ws.Range("C17").Value = 26
$ws.Range("C18").Value = 97
$ws.Range("C19").Value = 95
$ws.Range("C20").Value = 32
$ws.Range("C21").Value = 74

# "- loading" table (rows 23-42)
$ws.Range("C23").Value = 91
$ws.Range("C24").Value = 27
$ws.Range("C25").Value = 57
$ws.Range("C26").Value = 45
$ws.Range("C27").Value = 30
$ws.Range("C28").Value = 61
$ws.Range("C29").Value = 153
$ws.Range("C30").Value = 106
$ws.Range("C31").Value = 108
$ws.Range("C32").Value = 123
$ws.Range("C33").Value = 92
$ws.Range("C34").Value = 117
$ws.Range("C35").Value = 44
$ws.Range("C36").Value = 100
$ws.Range("C37").Value = 90
$ws.Range("C38").Value = 15
$ws.Range("C39").Value = 113
$ws.Range("C40").Value = 99
$ws.Range("C41").Value = 40
$ws.Range("C42").Value = 121

$wb.Save()
